$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 24.321797
$ws.Cells.Item(2, 8).Value = 72.96539100000001
$ws.Cells.Item(2, 9).Value = 0.005044792378505166
$ws.Cells.Item(2, 10).Value = 0.005044792378505166
$ws.Cells.Item(2, 13).Value = 2.270682
$ws.Cells.Item(2, 14).Value = 6.812046
$ws.Cells.Item(2, 15).Value = 0.1036812626940959
$ws.Cells.Item(2, 16).Value = 0.1036812626940959
$ws.Cells.Item(2, 17).Value = 55.227066655554
$ws.Cells.Item(2, 18).Value = 497.043599899986
$ws.Cells.Item(2, 19).Value = 0.0005230504438329669
$ws.Cells.Item(2, 20).Value = 0.0005230504438329669

$ws.Cells.Item(3, 7).Value = 24.321797
$ws.Cells.Item(3, 8).Value = 72.96539100000001
$ws.Cells.Item(3, 9).Value = 0.005044792378505166
$ws.Cells.Item(3, 10).Value = 0.005044792378505166
$ws.Cells.Item(3, 15).Value = 0.4406411276013061
$ws.Cells.Item(3, 16).Value = 0.4406411276013061
$ws.Cells.Item(3, 17).Value = 234.7127754126164
$ws.Cells.Item(3, 18).Value = 2112.414978713547
$ws.Cells.Item(3, 19).Value = 0.002222943002178991
$ws.Cells.Item(3, 20).Value = 0.002222943002178992

$ws.Cells.Item(4, 7).Value = 24.321797
$ws.Cells.Item(4, 8).Value = 72.96539100000001
$ws.Cells.Item(4, 9).Value = 0.005044792378505166
$ws.Cells.Item(4, 10).Value = 0.005044792378505166
$ws.Cells.Item(4, 13).Value = 9.928499666666665
$ws.Cells.Item(4, 14).Value = 29.785499
$ws.Cells.Item(4, 15).Value = 0.4533437011866523
$ws.Cells.Item(4, 16).Value = 0.4533437011866523
$ws.Cells.Item(4, 17).Value = 241.4789534072343
$ws.Cells.Item(4, 18).Value = 2173.310580665109
$ws.Cells.Item(4, 19).Value = 0.002287024848589747
$ws.Cells.Item(4, 20).Value = 0.002287024848589747

$ws.Cells.Item(5, 7).Value = 24.321797
$ws.Cells.Item(5, 8).Value = 72.96539100000001
$ws.Cells.Item(5, 9).Value = 0.005044792378505166
$ws.Cells.Item(5, 10).Value = 0.005044792378505166
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.051114
$ws.Cells.Item(5, 14).Value = 0.153342
$ws.Cells.Item(5, 15).Value = 0.002333908517945717
$ws.Cells.Item(5, 16).Value = 0.002333908517945718
$ws.Cells.Item(5, 17).Value = 1.243184331858
$ws.Cells.Item(5, 18).Value = 11.188658986722
$ws.Cells.Item(5, 19).Value = 0.00001177408390346084
$ws.Cells.Item(5, 20).Value = 0.00001177408390346084

$ws.Cells.Item(6, 9).Value = 0.971518017402211
$ws.Cells.Item(6, 10).Value = 0.971518017402211
$ws.Cells.Item(6, 13).Value = 2.270682
$ws.Cells.Item(6, 14).Value = 6.812046
$ws.Cells.Item(6, 15).Value = 0.1036812626940959
$ws.Cells.Item(6, 16).Value = 0.1036812626940959
$ws.Cells.Item(6, 17).Value = 10635.53983564373
$ws.Cells.Item(6, 18).Value = 95719.8585207936
$ws.Cells.Item(6, 19).Value = 0.1007282147743258
$ws.Cells.Item(6, 20).Value = 0.1007282147743259

$ws.Cells.Item(7, 9).Value = 0.971518017402211
$ws.Cells.Item(7, 10).Value = 0.971518017402211
$ws.Cells.Item(7, 15).Value = 0.4406411276013061
$ws.Cells.Item(7, 16).Value = 0.4406411276013061
$ws.Cells.Item(7, 19).Value = 0.4280907946730956
$ws.Cells.Item(7, 20).Value = 0.4280907946730956

$ws.Cells.Item(8, 9).Value = 0.971518017402211
$ws.Cells.Item(8, 10).Value = 0.971518017402211
$ws.Cells.Item(8, 13).Value = 9.928499666666665
$ws.Cells.Item(8, 14).Value = 29.785499
$ws.Cells.Item(8, 15).Value = 0.4533437011866523
$ws.Cells.Item(8, 16).Value = 0.4533437011866523
$ws.Cells.Item(8, 17).Value = 46503.62917969529
$ws.Cells.Item(8, 18).Value = 418532.6626172576
$ws.Cells.Item(8, 19).Value = 0.4404315737786368
$ws.Cells.Item(8, 20).Value = 0.4404315737786368

$ws.Cells.Item(9, 9).Value = 0.971518017402211
$ws.Cells.Item(9, 10).Value = 0.971518017402211
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.051114
$ws.Cells.Item(9, 14).Value = 0.153342
$ws.Cells.Item(9, 15).Value = 0.002333908517945717
$ws.Cells.Item(9, 16).Value = 0.002333908517945718
$ws.Cells.Item(9, 17).Value = 239.410442835718
$ws.Cells.Item(9, 18).Value = 2154.693985521462
$ws.Cells.Item(9, 19).Value = 0.002267434176152756
$ws.Cells.Item(9, 20).Value = 0.002267434176152756

$ws.Cells.Item(10, 7).Value = 108.9258753333333
$ws.Cells.Item(10, 8).Value = 326.777626
$ws.Cells.Item(10, 9).Value = 0.02259324940930984
$ws.Cells.Item(10, 10).Value = 0.02259324940930984
$ws.Cells.Item(10, 13).Value = 2.270682
$ws.Cells.Item(10, 14).Value = 6.812046
$ws.Cells.Item(10, 15).Value = 0.1036812626940959
$ws.Cells.Item(10, 16).Value = 0.1036812626940959
$ws.Cells.Item(10, 17).Value = 247.336024453644
$ws.Cells.Item(10, 18).Value = 2226.024220082796
$ws.Cells.Item(10, 19).Value = 0.002342496627119879
$ws.Cells.Item(10, 20).Value = 0.00234249662711988

$ws.Cells.Item(11, 7).Value = 108.9258753333333
$ws.Cells.Item(11, 8).Value = 326.777626
$ws.Cells.Item(11, 9).Value = 0.02259324940930984
$ws.Cells.Item(11, 10).Value = 0.02259324940930984
$ws.Cells.Item(11, 15).Value = 0.4406411276013061
$ws.Cells.Item(11, 16).Value = 0.4406411276013061
$ws.Cells.Item(11, 17).Value = 1051.167991975894
$ws.Cells.Item(11, 18).Value = 9460.511927783044
$ws.Cells.Item(11, 19).Value = 0.009955514895895828
$ws.Cells.Item(11, 20).Value = 0.00995551489589583

$ws.Cells.Item(12, 7).Value = 108.9258753333333
$ws.Cells.Item(12, 8).Value = 326.777626
$ws.Cells.Item(12, 9).Value = 0.02259324940930984
$ws.Cells.Item(12, 10).Value = 0.02259324940930984
$ws.Cells.Item(12, 13).Value = 9.928499666666665
$ws.Cells.Item(12, 14).Value = 29.785499
$ws.Cells.Item(12, 15).Value = 0.4533437011866523
$ws.Cells.Item(12, 16).Value = 0.4533437011866523
$ws.Cells.Item(12, 17).Value = 1081.470516938375
$ws.Cells.Item(12, 18).Value = 9733.234652445373
$ws.Cells.Item(12, 19).Value = 0.01024250730904967
$ws.Cells.Item(12, 20).Value = 0.01024250730904967

$ws.Cells.Item(13, 7).Value = 108.9258753333333
$ws.Cells.Item(13, 8).Value = 326.777626
$ws.Cells.Item(13, 9).Value = 0.02259324940930984
$ws.Cells.Item(13, 10).Value = 0.02259324940930984
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.051114
$ws.Cells.Item(13, 14).Value = 0.153342
$ws.Cells.Item(13, 15).Value = 0.002333908517945717
$ws.Cells.Item(13, 16).Value = 0.002333908517945718
$ws.Cells.Item(13, 17).Value = 5.567637191788
$ws.Cells.Item(13, 18).Value = 50.108734726092
$ws.Cells.Item(13, 19).Value = 0.00005273057724446026
$ws.Cells.Item(13, 20).Value = 0.00005273057724446028

$ws.Cells.Item(14, 7).Value = 4.068781333333334
$ws.Cells.Item(14, 8).Value = 12.206344
$ws.Cells.Item(14, 9).Value = 0.0008439408099740362
$ws.Cells.Item(14, 10).Value = 0.0008439408099740362
$ws.Cells.Item(14, 13).Value = 2.270682
$ws.Cells.Item(14, 14).Value = 6.812046
$ws.Cells.Item(14, 15).Value = 0.1036812626940959
$ws.Cells.Item(14, 16).Value = 0.1036812626940959
$ws.Cells.Item(14, 17).Value = 9.238908535536
$ws.Cells.Item(14, 18).Value = 83.15017681982401
$ws.Cells.Item(14, 19).Value = 0.00008750084881718609
$ws.Cells.Item(14, 20).Value = 0.00008750084881718611

$ws.Cells.Item(15, 7).Value = 4.068781333333334
$ws.Cells.Item(15, 8).Value = 12.206344
$ws.Cells.Item(15, 9).Value = 0.0008439408099740362
$ws.Cells.Item(15, 10).Value = 0.0008439408099740362
$ws.Cells.Item(15, 15).Value = 0.4406411276013061
$ws.Cells.Item(15, 16).Value = 0.4406411276013061
$ws.Cells.Item(15, 17).Value = 39.26498355749423
$ws.Cells.Item(15, 18).Value = 353.3848520174481
$ws.Cells.Item(15, 19).Value = 0.0003718750301357189
$ws.Cells.Item(15, 20).Value = 0.000371875030135719

$ws.Cells.Item(16, 7).Value = 4.068781333333334
$ws.Cells.Item(16, 8).Value = 12.206344
$ws.Cells.Item(16, 9).Value = 0.0008439408099740362
$ws.Cells.Item(16, 10).Value = 0.0008439408099740362
$ws.Cells.Item(16, 13).Value = 9.928499666666665
$ws.Cells.Item(16, 14).Value = 29.785499
$ws.Cells.Item(16, 15).Value = 0.4533437011866523
$ws.Cells.Item(16, 16).Value = 0.4533437011866523
$ws.Cells.Item(16, 17).Value = 40.39689411173956
$ws.Cells.Item(16, 18).Value = 363.572047005656
$ws.Cells.Item(16, 19).Value = 0.0003825952503760908
$ws.Cells.Item(16, 20).Value = 0.0003825952503760909

$ws.Cells.Item(17, 7).Value = 4.068781333333334
$ws.Cells.Item(17, 8).Value = 12.206344
$ws.Cells.Item(17, 9).Value = 0.0008439408099740362
$ws.Cells.Item(17, 10).Value = 0.0008439408099740362
$ws.Cells.Item(17, 11).Value = 2
$ws.Cells.Item(17, 12).Value = 0.6666666666666666
$ws.Cells.Item(17, 13).Value = 0.051114
$ws.Cells.Item(17, 14).Value = 0.153342
$ws.Cells.Item(17, 15).Value = 0.002333908517945717
$ws.Cells.Item(17, 16).Value = 0.002333908517945718
$ws.Cells.Item(17, 17).Value = 0.207971689072
$ws.Cells.Item(17, 18).Value = 1.871745201648
$ws.Cells.Item(17, 19).Value = 0.000001969680645040411
$ws.Cells.Item(17, 20).Value = 0.000001969680645040411
